# Applies the HTML-markup / content edits described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) D26 (HT203022 body): drop the stray leading "<br>" line before the
#        "<strong><h4>...</h4>" banner and bump the heading level h4 -> h1 ---
$d26 = $ws.Range("D26").Value2
$d26 = $d26.Replace("<br>`n<br><strong><h4>", "<br><strong><h1>")
$d26 = $d26.Replace("연식</h4>", "연식</h1>")
$ws.Range("D26").Value2 = $d26

# --- 2) D34 (HT203027 body): drop the leading "<br>" before "지원 모델" ---
$d34 = $ws.Range("D34").Value2
$d34 = $d34.Replace("<br>지원 모델: iCup 지원 차량", "지원 모델: iCup 지원 차량")
$ws.Range("D34").Value2 = $d34

# --- 3) D37 (warranty escalation body): wrap the final quoted <em> block in
#        <string> tags and move the closing quote to the end of the sentence ---
$d37 = $ws.Range("D37").Value2
$old3 = "<em>`"보증에 대한 확인을 진행하였으나, 동일하게 답변이 되는 점에 대해 사과드립니다.`n<br>해당 내용으로 추가적으로 확인했으나 보증 적용은 어렵습니다.</em>"
$new3 = "<em><string>`"보증에 대한 확인을 진행하였으나, 동일하게 답변이 되는 점에 대해 사과드립니다.`n<br>해당 내용으로 추가적으로 확인했으나 보증 적용은 어렵습니다.`"</em></string>"
$d37 = $d37.Replace($old3, $new3)
$ws.Range("D37").Value2 = $d37

# --- 4) Row-height adjustments that follow from the re-wrapped text above ---
$ws.Rows(24).RowHeight = 209.25
$ws.Rows(26).RowHeight = 244.5

# --- 5) New row 38: a single date value in column B, matching the other
#        entries dated 2024-06-11 (serial 45454) ---
$ws.Range("B38").Value2 = 45454

# --- 6) Move the active selection / viewport down to the newly added row ---
[void]$ws.Range("A38").Select()
